# "Update countries & provincias Spain"
#
# A later data pull for the Covid table (27 Sep 2020, 13:40 -> 14:57)
# refreshed several countries' case counts. Because the sheet is sorted
# by "Casos totales" descending, a handful of countries whose totals
# leap-frogged their neighbours end up keeping their row position while
# the country name + stats shift to the row above/below (i.e. the rows
# keep their numeric values, but the row each country now occupies
# changes). We therefore write both the country name (col A) and the
# stats (cols B-H) explicitly for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range("A1").Value2 = "Datos actualizados a 27 de Septiembre de 2020 a las 14:57"

# Estados Unidos -- new case counts
$ws.Range("B4").Value2 = 7288094
$ws.Range("C4").Value2 = 533
$ws.Range("D4").Value2 = 4524383
$ws.Range("E4").Value2 = 2554531
$ws.Range("G4").Value2 = 3
$ws.Range("H4").Value2 = 209180

# Catar -- new case counts
$ws.Range("B32").Value2 = 125084
$ws.Range("C32").Value2 = 234
$ws.Range("D32").Value2 = 121995
$ws.Range("E32").Value2 = 2875

# Ranking shuffle: Republica Dominicana / Panama / Paises Bajos swap rows 36-38
$ws.Range("A36").Value2 = "Paises Bajos"
$ws.Range("B36").Value2 = 111626
$ws.Range("C36").Value2 = 2995
$ws.Range("D36").Value2 = 0
$ws.Range("E36").Value2 = 0
$ws.Range("G36").Value2 = 8
$ws.Range("H36").Value2 = 6374

$ws.Range("A37").Value2 = "Republica Dominicana"
$ws.Range("B37").Value2 = 110957
$ws.Range("D37").Value2 = 85220
$ws.Range("E37").Value2 = 23644
$ws.Range("H37").Value2 = 2093

$ws.Range("A38").Value2 = "Panama"
$ws.Range("B38").Value2 = 110108
$ws.Range("D38").Value2 = 86796
$ws.Range("E38").Value2 = 20989
$ws.Range("H38").Value2 = 2323

# Kuwait -- new case counts
$ws.Range("B40").Value2 = 103544
$ws.Range("C40").Value2 = 345
$ws.Range("D40").Value2 = 94929
$ws.Range("E40").Value2 = 8014
$ws.Range("G40").Value2 = 4
$ws.Range("H40").Value2 = 601

# Bosnia y Herzegovina -- new case counts
$ws.Range("B79").Value2 = 26920
$ws.Range("C79").Value2 = 123
$ws.Range("E79").Value2 = 6352
$ws.Range("G79").Value2 = 2
$ws.Range("H79").Value2 = 822

# Dinamarca -- new case counts
$ws.Range("B80").Value2 = 26637
$ws.Range("C80").Value2 = 424
$ws.Range("D80").Value2 = 19650
$ws.Range("E80").Value2 = 6338
$ws.Range("G80").Value2 = 1
$ws.Range("H80").Value2 = 649

# Ranking shuffle: Gibraltar / Eritrea / San Martin (Parte Francesa) swap rows 180-182
$ws.Range("A180").Value2 = "San Martin (Parte Francesa)"
$ws.Range("B180").Value2 = 383
$ws.Range("C180").Value2 = 16
$ws.Range("D180").Value2 = 273
$ws.Range("E180").Value2 = 102
$ws.Range("H180").Value2 = 8

$ws.Range("A181").Value2 = "Gibraltar"
$ws.Range("B181").Value2 = 379
$ws.Range("C181").Value2 = 7
$ws.Range("E181").Value2 = 38

$ws.Range("A182").Value2 = "Eritrea"
$ws.Range("B182").Value2 = 375
$ws.Range("D182").Value2 = 341
$ws.Range("E182").Value2 = 34
$ws.Range("H182").Value2 = 0

# Ranking shuffle: Macao / San Bartolome swap rows 200-201
$ws.Range("A200").Value2 = "San Bartolome"
$ws.Range("B200").Value2 = 48
$ws.Range("C200").Value2 = 3
$ws.Range("D200").Value2 = 25
$ws.Range("E200").Value2 = 23

$ws.Range("A201").Value2 = "Macao"
$ws.Range("B201").Value2 = 46
$ws.Range("D201").Value2 = 46
$ws.Range("E201").Value2 = 0

# Ranking shuffle: Nueva Caledonia / Timor Oriental / Dominica swap rows 205,206,208
$ws.Range("A205").Value2 = "Dominica"
$ws.Range("B205").Value2 = 30
$ws.Range("C205").Value2 = 6
$ws.Range("D205").Value2 = 24
$ws.Range("E205").Value2 = 6

$ws.Range("A206").Value2 = "Nueva Caledonia"
$ws.Range("D206").Value2 = 26
$ws.Range("E206").Value2 = 1

$ws.Range("A208").Value2 = "Timor Oriental"
$ws.Range("B208").Value2 = 27
$ws.Range("D208").Value2 = 27
$ws.Range("E208").Value2 = 0
